$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$app = $excel

# Switch to manual calculation while we update the raw inputs so that
# dependent formulas (rows 4/6/8, 11, 28, ...) are only recalculated once,
# against the final, fully-updated input values - this matches how the
# pollster bias bug fix replaced the whole block of primary-vote inputs.
$app.Calculation = -4135  # xlCalculationManual

# --- Row 2: raw poll primary votes for this poll (bug fix: correct party
#     columns now receive numbers, previously-numeric columns now correctly
#     show #N/A where that pollster doesn't report a breakdown) ---
$ws.Range("A2").Value = 35
$ws.Range("B2").Value = 30
$ws.Range("C2").Value = 13
$ws.Range("D2").Value = 7
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = "#N/A"
$ws.Range("H2").Value = "#N/A"
$ws.Range("I2").Value = 8

# --- Row 10: two-party-preferred raw inputs ---
$ws.Range("A10").Value = 48
$ws.Range("B10").Value = 47

# --- New section: pollster bias / preference-flow survival check ---
$ws.Range("A22").Value = "OPV"

$ws.Range("A23").Value = "LNP"
$ws.Range("B23").Value = "ALP"
$ws.Range("C23").Value = "GRN"
$ws.Range("D23").Value = "ONP"
$ws.Range("E23").Value = "NXT"
$ws.Range("F23").Value = "UAP"
$ws.Range("G23").Value = "DEM"
$ws.Range("H23").Value = "DLP"
$ws.Range("I23").Value = "OTH"

$ws.Range("A24").Value = 41.6
$ws.Range("B24").Value = 42.25
$ws.Range("C24").Value = 8.37
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 7.78

$ws.Range("A25").Value = 100
$ws.Range("B25").Value = 100
$ws.Range("C25").Value = 60
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 60

$ws.Range("A26").Value = 0
$ws.Range("B26").Value = 100
$ws.Range("C26").Value = 70
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 40

# Add the two label strings in the same order the original author typed
# them so the shared-string table indices line up ("<- pref flow" before
# "<- survival").
$ws.Range("J26").Value = "<- pref flow"
$ws.Range("J25").Value = "<- survival"

$ws.Range("I28").Formula = "=SUMPRODUCT(A24:I24,A25:I25,A26:I26)/(SUMPRODUCT(A24:I24,A25:I25))"

# Restore automatic calculation and force a full recalculation so every
# cached formula value (rows 4, 6, 8, 11, 28, ...) reflects the corrected
# inputs above.
$app.Calculation = -4105  # xlCalculationAutomatic
$app.CalculateFullRebuild()

# Leave the selection where the author left it after making this edit.
$ws.Range("I27").Select()
